# Applies the commit "Ajout du CR, des graphiques, des etats de Art" to
# Presentation_Objectifs_Rendus_20sept.pptx:
#   - Slide 5 ("Ensemble de données"), "Content Placeholder 7":
#       splits the "datasets" run out of paragraph 1, and the "GIFs" run
#       out of paragraph 2 (so spell-check can flag them individually).
#   - Slide 6, "Content Placeholder 2": the text is unchanged (only
#     PowerPoint-internal proofing flags changed in the source diff).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 - "Content Placeholder 7"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(3)
$tr5 = $shp5.TextFrame.TextRange

# Paragraph 1: "...En considérant ces données comme des datasets : "
# -> split into 3 runs: "...comme des ", "datasets", "<nbsp>: "
$para1 = $tr5.Paragraphs(1, 1)
$datasetsRun = $para1.Characters(265, 8)
$datasetsRun.Text = $datasetsRun.Text

# Paragraph 2: "Liste / Catégories / GIFs des exercices"
# -> split into 3 runs: "Liste / Catégories / ", "GIFs", " des exercices"
$para2 = $tr5.Paragraphs(2, 1)
$gifsRun = $para2.Characters(22, 4)
$gifsRun.Text = $gifsRun.Text
